# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" stats sheet: new totals for several countries
# and a re-sort that swaps the ranking of three country pairs (Zambia/
# Eslovenia, Belice/Lesoto, Laos/Santa Lucia), plus the "last updated" stamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner (row 1)
$ws.Range("A1").Value = "Datos actualizados a 1 de Julio de 2020 a las 16:29"

# Helper: write one data row (country name + 7 numeric columns B..H)
function Set-Row([int]$row, [string]$country, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $country
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

Set-Row 4   "Estados Unidos" 2731671 3818 1143688 1457826 0 35 130157
Set-Row 5   "Brasil"         1409693 1208 790040  559908  0 89 59745
Set-Row 7   "India"          591420  5628 352171  221754  0 85 17495
Set-Row 24  "Catar"          97003   915  83965   12923   0 2  115
Set-Row 28  "Argentina"      64530   0    23040   40180   0 3  1310
Set-Row 62  "Serbia"         14836   272  12772   1783    0 4  281
Set-Row 76  "Finlandia"      7236    22   6700    208     0 0  328
Set-Row 80  "Kenia"          6673    307  2089    4435    0 1  149
Set-Row 111 "Sri Lanka"      2050    3    1748    291     0 0  11
Set-Row 115 "Islandia"       1847    5    1823    14      0 0  10
Set-Row 117 "Libano"         1788    10   1223    531     0 0  34

# Zambia / Eslovenia swap ranking (row 120 <-> 121), with fresh Zambia data
Set-Row 120 "Zambia"         1632    38   1348    254     0 6  30
Set-Row 121 "Eslovenia"      1613    13   1384    118     0 0  111

Set-Row 154 "Reunion"        528     2    472     54      0 0  2
Set-Row 155 "Surinam"        517     2    227     277     0 0  13
Set-Row 162 "Birmania"       303     4    222     75      0 0  6

# Belice / Lesoto swap ranking (row 197 <-> 198), with fresh Belice data
Set-Row 197 "Belice"         28      4    18      8       0 0  2
Set-Row 198 "Lesoto"         27      0    4       23      0 0  0

# Laos / Santa Lucia swap ranking (row 203 <-> 204); identical stats either way
Set-Row 203 "Laos"           19      0    19      0       0 0  0
Set-Row 204 "Santa Lucia"    19      0    19      0       0 0  0
